$wb = $excel.ActiveWorkbook
for ($i = 1; $i -le $wb.Worksheets.Count; $i++) {
    $s = $wb.Worksheets.Item($i)
    Write-Host "Index $i -> $($s.Index) dim"
}
